$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 32 (pushes the existing rows 32-35 down to 33-36) ---
$ws.Range("A32").EntireRow.Insert()

# Re-use formatting from the row that is now just below the new blank row
# (the old row 32, shifted down to row 33) so the new cells keep the same
# border/format family rather than falling back to "no style".
$ws.Range("A33:E33").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the new DRAIAM103 test case content ---
$ws.Range("A32").Value = "DRAIAM103"
$ws.Range("B32").Value = "OPQA-5139 || OPQA-5140"
$ws.Range("C32").Value = "Verify that 'Project Neon' should be moved within the white area and should be above 'Forgot Password' text and center aligned`n|| Verify that Clarivate Analytics logo should be below the Marketing area."
$ws.Range("D32").Value = "Y"

# The row holds wrapped, two-line text, same as similar rows - give it the
# matching row height.
$ws.Rows(32).RowHeight = 30

# --- Update the view: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("C34").Select()
